$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shift existing data one column to the right (A->B, B->C) and fold
# --- the old duplicate D column into C, then insert the new Timestamp
# --- column in A together with its header.
# (.Value2 is used throughout -- the plain .Value getter mis-behaves on
# this host for numeric cells.)

# Header row
$ws.Cells.Item(1, 4).Value2 = $ws.Cells.Item(1, 3).Value2   # D1 = old C1 ("variable 3")
$ws.Cells.Item(1, 3).Value2 = $ws.Cells.Item(1, 2).Value2   # C1 = old B1 ("variable 2")
$ws.Cells.Item(1, 2).Value2 = $ws.Cells.Item(1, 1).Value2   # B1 = old A1 ("variable 1")
$ws.Cells.Item(1, 1).Value2 = "Timestamp"

# Data rows 2-4: shift B<-A, C<-B, D<-C (old D was a duplicate of old C, now dropped)
for ($r = 2; $r -le 4; $r++) {
    $ws.Cells.Item($r, 4).Value2 = $ws.Cells.Item($r, 3).Value2
    $ws.Cells.Item($r, 3).Value2 = $ws.Cells.Item($r, 2).Value2
    $ws.Cells.Item($r, 2).Value2 = $ws.Cells.Item($r, 1).Value2
}

# New timestamp values in column A, stored as text
$ws.Cells.Item(2, 1).Value2 = "08.02.2018 13:34:20"
$ws.Cells.Item(3, 1).Value2 = "08.02.2018 13:34:30"
$ws.Cells.Item(4, 1).Value2 = "08.02.2018 13:34:40"
$ws.Range("A2:A4").NumberFormat = "@"

# Drop the now-obsolete 5th row entirely
$ws.Rows("5:5").Delete()

# Column widths: column D now best-fits the numeric values instead of the
# fixed width it shared with B:C (12.5703125, i.e. defaultColWidth + 1)
$ws.Columns("D:D").ColumnWidth = 11.736979166666666

# Selection follows the last data row as in the authored workbook
$ws.Range("A4").Select() | Out-Null
